$d = $word.ActiveDocument

# Move to the very end of the document body (after the last paragraph's text)
# and insert a brand-new paragraph there containing the new sentence.
$end = $d.Content
$end.Collapse(0) | Out-Null

$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>埃斯发士大夫fd</w:t></w:r></w:p>'

$end.InsertXML($newParaXml) | Out-Null
